# Apply "Automatic update of files" changes to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Förändrad date bump + updated counts + new species line.
$ws.Range("C2").Value = 45179
$ws.Range("K2").Value = 2
$ws.Range("O2").Value = 11
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 16

$r2Original = $ws.Range("R2").Value2
$ws.Range("R2").Value = "Grangråticka`r`n" + $r2Original

# Rows 3-15: only the "Förändrad" date (column C) moves forward one day.
$ws.Range("C3").Value = 45179
$ws.Range("C4").Value = 45179
$ws.Range("C5").Value = 45179
$ws.Range("C6").Value = 45179
$ws.Range("C7").Value = 45179
$ws.Range("C8").Value = 45179
$ws.Range("C9").Value = 45179
$ws.Range("C10").Value = 45179
$ws.Range("C11").Value = 45179
$ws.Range("C12").Value = 45179
$ws.Range("C13").Value = 45179
$ws.Range("C14").Value = 45179
$ws.Range("C15").Value = 45179
